# grid_template.xlsx maintenance edit ("Miglioramenti e pulizia generale")
#  1. A3's placeholder text changes from ${i} to ${String}
#  2. A new cell comment is added on A3 documenting the jx:each() tags that
#     replace the old jx:area() comment's single jx:area() tag
#  3. The sheet's active selection moves from A4 to G8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the A3 placeholder text
$ws.Range("A3").Value = '${String}'

# 2. Add the explanatory comment to A3 (mirrors the existing A1 comment style)
$newline = [char]10
$line1 = 'jx:each(items="master" var="items" lastCell="A3" direction="RIGHT")'
$line2 = 'jx:each(items="items" var="String" lastCell="A3" direction="DOWN")'
$commentText = "Author:" + $newline + $line1 + $newline + $line2
$comment = $ws.Range("A3").AddComment($commentText)

# Best-effort: bold the "Author:" label, matching the A1 comment's formatting
$comment.Shape.TextFrame.Characters(1, 7).Font.Bold = $true

# 3. Move the active selection to G8
$ws.Range("G8").Select()
